$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume data refresh (GitHub Actions scheduled update).
# D-column cells hold price text that can look numeric ("2.43", "61.34", ...);
# force Text format before assigning so Excel does not silently convert them
# to numbers (matches the source data which stores these as plain strings).
$priceCells = @("D2", "D3", "D5", "D7", "D9", "D10", "D12", "D13", "D14", "D15", "D17", "D18", "D19", "D21", "D22", "D24", "D26", "D27", "D28", "D29", "D32", "D33", "D34", "D36", "D40", "D41", "D42", "D45", "D46", "D47", "D49", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "38.610.30"
$ws.Range("E2").Value = "  +2.39%  "
$ws.Range("D3").Value = "2.091.90"
$ws.Range("E3").Value = "  +2.95%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "228.90"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("E6").Value = "  +1.11%  "
$ws.Range("D7").Value = "61.34"
$ws.Range("E7").Value = "  +2.16%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.381"
$ws.Range("D10").Value = "0.0842"
$ws.Range("E10").Value = "  +2.39%  "
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").Value = "2.399.15"
$ws.Range("E12").Value = "  +2.80%  "
$ws.Range("D13").Value = "14.79"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "22.31"
$ws.Range("E14").Value = "  +6.12%  "
$ws.Range("D15").Value = "0.782"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("E16").Value = "  +5.22%  "
$ws.Range("D17").Value = "2.105.25"
$ws.Range("E17").Value = "  +3.71%  "
$ws.Range("D18").Value = "38.563.39"
$ws.Range("E18").Value = "  +2.35%  "
$ws.Range("D19").Value = "71.01"
$ws.Range("E19").Value = "  +2.09%  "
$ws.Range("E20").Value = "  +3.21%  "
$ws.Range("D21").Value = "0.0₃0834"
$ws.Range("E21").Value = "  +1.43%  "
$ws.Range("D22").Value = "226.04"
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("D26").Value = "170.63"
$ws.Range("E26").Value = "  +1.97%  "
$ws.Range("D27").Value = "9.44"
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("D28").Value = "0.134"
$ws.Range("E28").Value = "  +3.24%  "
$ws.Range("D29").Value = "19.09"
$ws.Range("E29").Value = "  +1.80%  "
$ws.Range("E30").Value = "  +7.61%  "
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("D32").Value = "2.27"
$ws.Range("E32").Value = "  +2.54%  "
$ws.Range("D33").Value = "4.80"
$ws.Range("E33").Value = "  +6.83%  "
$ws.Range("D34").Value = "4.49"
$ws.Range("E34").Value = "  +2.70%  "
$ws.Range("E35").Value = "  +0.76%  "
$ws.Range("D36").Value = "6.54"
$ws.Range("E36").Value = "  +1.74%  "
$ws.Range("E37").Value = "  +3.01%  "
$ws.Range("E38").Value = "  +4.78%  "
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").Value = "18.55"
$ws.Range("E40").Value = "  +2.59%  "
$ws.Range("D41").Value = "1.545.69"
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("D42").Value = "99.79"
$ws.Range("E42").Value = "  +4.04%  "
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("E44").Value = "  +1.20%  "
$ws.Range("D45").Value = "0.0910"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "7.69"
$ws.Range("E46").Value = "  +9.37%  "
$ws.Range("D47").Value = "4.14"
$ws.Range("E47").Value = "  +3.04%  "
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("D49").Value = "1.03"
$ws.Range("E49").Value = "  +2.89%  "
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("D51").Value = "2.290.97"
$ws.Range("E51").Value = "  +3.07%  "
